# Apply numeric updates to the "F" column (follower/sales count) cells
# across the four worksheets, as described by the source diff.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 396
$ws1.Range("F5").Value = 1345
$ws1.Range("F7").Value = 2562
$ws1.Range("F9").Value = 18924
$ws1.Range("F11").Value = 2002

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 115

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 592
$ws3.Range("F4").Value = 565

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 592
$ws4.Range("F5").Value = 565
$ws4.Range("F9").Value = 1345
$ws4.Range("F14").Value = 2562
$ws4.Range("F16").Value = 18924
$ws4.Range("F22").Value = 2002
$ws4.Range("F24").Value = 115
